$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.513.08"
$ws.Range("E2").Value = "  -0.62%  "

$ws.Range("D3").Value = "2.579.29"
$ws.Range("E3").Value = "  -1.70%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.36%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.78%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("E8").Value = "  -1.45%  "

$ws.Range("D9").Value = "2.578.50"
$ws.Range("E9").Value = "  -1.66%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.136"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.99%  "

$ws.Range("E11").Value = "  +0.29%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.357"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.48%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.19"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.60%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.73"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.45%  "

$ws.Range("D15").Value = "3.056.53"
$ws.Range("E15").Value = "  -1.41%  "

$ws.Range("E16").Value = "  -1.94%  "

$ws.Range("D17").Value = "66.289.85"
$ws.Range("E17").Value = "  -0.79%  "

$ws.Range("D18").Value = "2.581.49"
$ws.Range("E18").Value = "  -1.61%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.41"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.09%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.45%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "352.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.60%  "

$ws.Range("E22").Value = "  -2.14%  "

$ws.Range("E23").Value = "  -0.85%  "

$ws.Range("E24").Value = "  +0.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.88"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.99%  "

$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "68.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.87%  "

$ws.Range("B27").Value = "Aptos"
$ws.Range("C27").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -8.23%  "

$ws.Range("D28").Value = "2.711.81"
$ws.Range("E28").Value = "  -1.59%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.13%  "

$ws.Range("D30").Value = "0.0₃0987"
$ws.Range("E30").Value = "  -2.13%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "535.98"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.58%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.03"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.66%  "

$ws.Range("E33").Value = "  -2.06%  "

$ws.Range("E34").Value = "  -2.06%  "

$ws.Range("E35").Value = "  -1.24%  "

$ws.Range("E36").Value = "  +0.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.46"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.94%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "156.75"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.40%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.74"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.04%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.360"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.54%  "

$ws.Range("E41").Value = "  +1.89%  "

$ws.Range("E42").Value = "  -0.05%  "

$ws.Range("E43").Value = "  -1.22%  "

$ws.Range("E44").Value = "  +0.06%  "

$ws.Range("E45").Value = "  -1.24%  "

$ws.Range("D46").Value = "0.0₆0286"
$ws.Range("E46").Value = "  -3.78%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "149.01"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.79%  "

$ws.Range("E48").Value = "  -2.95%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.72"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.04%  "

$ws.Range("E50").Value = "  -0.84%  "

$ws.Range("E51").Value = "  -1.38%  "
